$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.337.22"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.365.68"
$ws.Range("E3").Value = "  +4.58%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.25"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.648"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.01"
$ws.Range("E7").Value = "  +12.88%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.478"
$ws.Range("E9").Value = "  +6.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0982"
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "27.16"
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.724.36"
$ws.Range("E12").Value = "  +4.86%  "
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.98"
$ws.Range("E14").Value = "  +2.21%  "
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.861"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.360.99"
$ws.Range("E17").Value = "  +4.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.266.93"
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000102"
$ws.Range("E19").Value = "  +4.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.35"
$ws.Range("E20").Value = "  +3.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.34"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "249.76"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.99"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.42"
$ws.Range("E28").Value = "  +2.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.72"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E30").Value = "  +5.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.130"
$ws.Range("E31").Value = "  -5.21%  "
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.96"
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0692"
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("E36").Value = "  +6.30%  "
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.57"
$ws.Range("E37").Value = "  +1.64%  "
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.90"
$ws.Range("E41").Value = "  +2.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.48"
$ws.Range("E42").Value = "  +7.06%  "
$ws.Range("E43").Value = "  +8.38%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.98"
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.51"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0953"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.447.52"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.591.75"
$ws.Range("E49").Value = "  +4.78%  "
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("E51").Value = "  +0.41%  "
